$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target values look numeric but must remain exact text
# (Excel would otherwise normalize e.g. "193.40" -> 193.4). Temporarily mark as Text.
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"

$ws.Range("D2").Value = '26.108.35'
$ws.Range("E2").Value = '  +1.41%  '

$ws.Range("D3").Value = '1.644.73'
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("D4").Value = '0.994'
$ws.Range("E4").Value = '  -0.73%  '

$ws.Range("D5").Value = '215.97'
$ws.Range("E5").Value = '  +0.54%  '

$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.55%  '

$ws.Range("D8").Value = '0.258'
$ws.Range("E8").Value = '  -0.26%  '

$ws.Range("D9").Value = '0.0634'
$ws.Range("E9").Value = '  -0.38%  '

$ws.Range("D10").Value = '19.85'
$ws.Range("E10").Value = '  +0.85%  '

$ws.Range("D11").Value = '0.0790'
$ws.Range("E11").Value = '  +0.45%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '4.26'
$ws.Range("E12").Value = '  +0.44%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.660.32'
$ws.Range("E13").Value = '  +1.46%  '

$ws.Range("D14").Value = '1.866.03'
$ws.Range("E14").Value = '  +0.46%  '

$ws.Range("D15").Value = '0.555'
$ws.Range("E15").Value = '  -0.52%  '

$ws.Range("D16").Value = '0.0₃0764'
$ws.Range("E16").Value = '  +0.12%  '

$ws.Range("D17").Value = '63.34'
$ws.Range("E17").Value = '  +1.05%  '

$ws.Range("D18").Value = '26.115.12'
$ws.Range("E18").Value = '  +1.40%  '

$ws.Range("D19").Value = '0.998'
$ws.Range("E19").Value = '  -0.37%  '

$ws.Range("D20").Value = '4.47'
$ws.Range("E20").Value = '  +0.58%  '

$ws.Range("D21").Value = '193.40'
$ws.Range("E21").Value = '  -0.31%  '

$ws.Range("D22").Value = '10.02'
$ws.Range("E22").Value = '  +0.83%  '

$ws.Range("D23").Value = '6.40'
$ws.Range("E23").Value = '  +2.20%  '

$ws.Range("D24").Value = '0.996'
$ws.Range("E24").Value = '  -0.58%  '

$ws.Range("E25").Value = '  -1.52%  '

$ws.Range("D26").Value = '142.30'
$ws.Range("E26").Value = '  -0.29%  '

$ws.Range("E27").Value = '  +0.89%  '

$ws.Range("D28").Value = '6.90'
$ws.Range("E28").Value = '  +0.37%  '

$ws.Range("D29").Value = '15.64'
$ws.Range("E29").Value = '  +0.62%  '

$ws.Range("E30").Value = '  +0.32%  '

$ws.Range("D31").Value = '0.0496'
$ws.Range("E31").Value = '  +0.58%  '

$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  +0.69%  '

$ws.Range("E33").Value = '  +0.42%  '

$ws.Range("E34").Value = '  +1.60%  '

$ws.Range("D35").Value = '2.39'
$ws.Range("E35").Value = '  +0.17%  '

$ws.Range("D36").Value = '0.912'
$ws.Range("E36").Value = '  +1.22%  '

$ws.Range("D37").Value = '1.149.12'
$ws.Range("E37").Value = '  +1.94%  '

$ws.Range("D38").Value = '0.548'
$ws.Range("E38").Value = '  +0.41%  '

$ws.Range("E39").Value = '  -1.11%  '

$ws.Range("E40").Value = '  +0.86%  '

$ws.Range("D41").Value = '0.997'
$ws.Range("E41").Value = '  -0.47%  '

$ws.Range("E42").Value = '  +1.94%  '

$ws.Range("D43").Value = '101.00'
$ws.Range("E43").Value = '  +1.21%  '

$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("D45").Value = '1.776.66'
$ws.Range("E45").Value = '  +0.51%  '

$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₆0109'
$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '55.72'
$ws.Range("E47").Value = '  +1.23%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '1.47'
$ws.Range("E48").Value = '  +6.43%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.0512'
$ws.Range("E49").Value = '  +1.85%  '

$ws.Range("D50").Value = '7.68'
$ws.Range("E50").Value = '  +1.80%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.417'
$ws.Range("E51").Value = '  +0.05%  '

# Restore default (General) style for the forced-text cells so no stray number format remains
$ws.Range("D11").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D43").Style = "Normal"
